# Commit: "Fruta / hortaliza, semanal"
# A new daily price record was inserted as the new row 5 of the sheet,
# pushing the previously existing rows 5..61 down to rows 6..62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5 (shifts rows 5..61 -> 6..62,
# carrying along their formatting, and extends the used range to row 62).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = "2023-04-27"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107001
$ws.Range("J5").Value = "Caqui"
$ws.Range("K5").Value = "Fuyu"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 24000
$ws.Range("Q5").Value = "$/bandeja 15 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1600
$ws.Range("T5").Value = 15
